$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 90.40000000000001
$ws.Range("I4").Value = 88
$ws.Range("K4").Value = 88
$ws.Range("M4").Value = 26

$ws.Range("H12").Value = 148.875
$ws.Range("I12").Value = 89.75
$ws.Range("J12").Value = 208
$ws.Range("K12").Value = 89.75
$ws.Range("L12").Value = 208
$ws.Range("M12").Value = 80.25
$ws.Range("N12").Value = -548

$ws.Range("H40").Value = 2071.7144
$ws.Range("I40").Value = 1750
$ws.Range("J40").Value = 2200.4
$ws.Range("K40").Value = 1750
$ws.Range("L40").Value = 2200.4
$ws.Range("M40").Value = -1575
$ws.Range("N40").Value = -2550.4

$ws.Range("H125").Value = 2155.6155
$ws.Range("I125").Value = 4377.4
$ws.Range("J125").Value = 767
$ws.Range("K125").Value = 39396.6
$ws.Range("L125").Value = 6903
$ws.Range("M125").Value = -36936.6
$ws.Range("N125").Value = -11823

$ws.Range("H127").Value = 896.1539
$ws.Range("I127").Value = 766.6667
$ws.Range("J127").Value = 913.04346
$ws.Range("K127").Value = 2300.0001
$ws.Range("L127").Value = 2739.13038
$ws.Range("M127").Value = 2659.9999
$ws.Range("N127").Value = -12659.13038

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 26320562
$ws.Range("I61").Value = 38466210
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 38466210
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -38465998
$ws.Range("N61").Value = -5424

$ws.Range("H97").Value = 6243.8335
$ws.Range("I97").Value = 8451.23
$ws.Range("J97").Value = 504.6
$ws.Range("K97").Value = 8451.23
$ws.Range("L97").Value = 504.6
$ws.Range("M97").Value = -7955.23
$ws.Range("N97").Value = -1496.6

$ws.Range("H132").Value = 6412333.5
$ws.Range("I132").Value = 11365480
$ws.Range("J132").Value = 2379.5293
$ws.Range("K132").Value = 34096440
$ws.Range("L132").Value = 7138.5879
$ws.Range("M132").Value = -34093910
$ws.Range("N132").Value = -12198.5879

$ws.Range("H136").Value = 26320562
$ws.Range("I136").Value = 38466210
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 115398630
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -115396080
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1099.1666
$ws.Range("I99").Value = 898.75
$ws.Range("K99").Value = 898.75
$ws.Range("M99").Value = 599.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 96.5
$ws.Range("I7").Value = 83.42856999999999
$ws.Range("K7").Value = 83.42856999999999
$ws.Range("M7").Value = 29.57143000000001

$ws.Range("H88").Value = 47910.285
$ws.Range("J88").Value = 47910.285
$ws.Range("L88").Value = 47910.285
$ws.Range("N88").Value = -48722.285

$ws.Range("H91").Value = 47910.285
$ws.Range("J91").Value = 47910.285
$ws.Range("L91").Value = 47910.285
$ws.Range("N91").Value = -50718.285

$ws.Range("H99").Value = 1466.3334
$ws.Range("I99").Value = 1299.5
$ws.Range("J99").Value = 1800
$ws.Range("K99").Value = 1299.5
$ws.Range("L99").Value = 1800
$ws.Range("M99").Value = 198.5
$ws.Range("N99").Value = -4796

$ws.Range("H126").Value = 1466.3334
$ws.Range("I126").Value = 1299.5
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 3898.5
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -1428.5
$ws.Range("N126").Value = -10340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 145
$ws.Range("I2").Value = 187.5
$ws.Range("J2").Value = 81.25
$ws.Range("K2").Value = 1125
$ws.Range("L2").Value = 487.5
$ws.Range("M2").Value = -1012
$ws.Range("N2").Value = -713.5

$ws.Range("H12").Value = 35.6
$ws.Range("I12").Value = 15.2
$ws.Range("J12").Value = 42.4
$ws.Range("K12").Value = 45.59999999999999
$ws.Range("L12").Value = 127.2
$ws.Range("M12").Value = 127.4
$ws.Range("N12").Value = -473.2

$ws.Range("H17").Value = 740.2
$ws.Range("I17").Value = 100.5
$ws.Range("J17").Value = 1166.6666
$ws.Range("K17").Value = 301.5
$ws.Range("L17").Value = 3499.9998
$ws.Range("M17").Value = -132.5
$ws.Range("N17").Value = -3837.9998

$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H131").Value = 882.0606
$ws.Range("I131").Value = 895.55554
$ws.Range("J131").Value = 879.9298
$ws.Range("K131").Value = 2686.66662
$ws.Range("L131").Value = 2639.7894
$ws.Range("M131").Value = 2353.33338
$ws.Range("N131").Value = -12719.7894

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 27940
$ws.Range("J51").Value = 27940
$ws.Range("L51").Value = 27940
$ws.Range("N51").Value = -28958

$ws.Range("H97").Value = 1200
$ws.Range("I97").Value = 1248.3334
$ws.Range("K97").Value = 1248.3334
$ws.Range("M97").Value = -752.3334

$ws.Range("H141").Value = 269619.34
$ws.Range("J141").Value = 269619.34
$ws.Range("L141").Value = 269619.34
$ws.Range("N141").Value = -279979.34

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3074.0227
$ws.Range("I40").Value = 5045.643
$ws.Range("J40").Value = 2153.9333
$ws.Range("K40").Value = 5045.643
$ws.Range("L40").Value = 2153.9333
$ws.Range("M40").Value = -4909.643
$ws.Range("N40").Value = -2425.9333

$ws.Range("H100").Value = 2100
$ws.Range("I100").Value = 2000
$ws.Range("K100").Value = 2000
$ws.Range("M100").Value = -1459

$ws.Range("H122").Value = 8030.2104
$ws.Range("I122").Value = 10471.75
$ws.Range("J122").Value = 6254.5454
$ws.Range("K122").Value = 31415.25
$ws.Range("L122").Value = 18763.6362
$ws.Range("M122").Value = -28965.25
$ws.Range("N122").Value = -23663.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 889.6
$ws.Range("I81").Value = 967
$ws.Range("J81").Value = 580
$ws.Range("K81").Value = 1934
$ws.Range("L81").Value = 1160
$ws.Range("M81").Value = -873
$ws.Range("N81").Value = -3282

$ws.Range("H84").Value = 889.6
$ws.Range("I84").Value = 967
$ws.Range("J84").Value = 580
$ws.Range("K84").Value = 9670
$ws.Range("L84").Value = 5800
$ws.Range("M84").Value = -4366
$ws.Range("N84").Value = -16408
